# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
# 1. Insert a new "Player Info" worksheet as the first sheet, with the
#    player's ID / NAME / BATTING_HAND / BOWL_STYLE.
# 2. Rename the MATCH_CARD_LINK column to MATCH_CODE on both "ODI Batting"
#    and "ODI Bowling", replacing the full scorecard URL with just the
#    numeric match code that used to be the query-string parameter.

$wb = $excel.ActiveWorkbook

# --- 1. New "Player Info" sheet, inserted before "ODI Batting" ---------
$battingForInsert = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingForInsert)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "5851"
$playerInfo.Range("B2").Value = "Ravi Bishnoi"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Leg Break"

$playerInfo.Range("A1").Select()

# Re-fetch the other sheets by name now that the sheet order has changed.
$batting = $wb.Worksheets.Item("ODI Batting")
$bowling = $wb.Worksheets.Item("ODI Bowling")

# --- 2. ODI Batting: MATCH_CARD_LINK -> MATCH_CODE ----------------------
$batting.Range("D1").Value = "MATCH_CODE"
$batting.Range("D2").NumberFormat = "@"
$batting.Range("D2").Value = "4656"

# --- 3. ODI Bowling: MATCH_CARD_LINK -> MATCH_CODE ----------------------
$bowling.Range("B1").Value = "MATCH_CODE"
$bowling.Range("B2").NumberFormat = "@"
$bowling.Range("B2").Value = "4656"
